$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New QR-scanner log rows to append (Student ID, Subject, Log Date, Log Time, Type, User).
# All values must land in the sheet as literal text (matching the existing rows), so we
# avoid $range.Value = "..." (which lets Excel auto-coerce numeric-/date-/time-looking
# strings into numbers/dates) and instead write each value as a formula that evaluates to
# a string, then copy/paste-special the cached value back onto the target cell. That keeps
# the text as text without leaving any formula behind.

$newRows = @(
    @("234928", "Anatomy", "05/11/2025", "09:07:54", "Manual", "nahla.nagiub@med.asu.edu.eg"),
    @("234089", "Anatomy", "05/11/2025", "09:17:39", "Scan",   "nahla.nagiub@med.asu.edu.eg"),
    @("234137", "Anatomy", "05/11/2025", "09:19:56", "Scan",   "nahla.nagiub@med.asu.edu.eg")
)

$startRow = 82
$scratch = $ws.Range("Z1")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]

    for ($col = 1; $col -le 6; $col++) {
        $text = $rowValues[$col - 1]
        $target = $ws.Cells.Item($rowIndex, $col)

        $scratch.Formula = '="' + $text + '"'
        $scratch.Copy()
        $target.PasteSpecial(-4163)
    }
}

$scratch.ClearContents()
